# Update smoke test data for QA in Maestro data workbook

$wb = $excel.ActiveWorkbook

$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsHogar  = $wb.Worksheets.Item("DatosHogar")
$wsMotor  = $wb.Worksheets.Item("DatosMotor")
$wsAP     = $wb.Worksheets.Item("DatosAP")

# DatosCuenta sheet: Name / LastName / Documento / NumeroCalle
$wsCuenta.Range("A2").Value = "SmokQACuatro"
$wsCuenta.Range("B2").Value = "SmokeNameQACUatro"
$wsCuenta.Range("C2").Value = 27100120
$wsCuenta.Range("D2").Value = 122

# DatosHogar sheet: NvoNro
$wsHogar.Range("A2").Value = 641

# DatosMotor sheet: Patente / Motor / Chasis
$wsMotor.Range("A2").Value = "SMP022"
$wsMotor.Range("B2").Value = "ABC12SSMP022"
$wsMotor.Range("C2").Value = "ZAZ123SSMP022"

# DatosAP sheet: NumeroDocumento
$wsAP.Range("A2").Value = 21200123

# Make DatosCuenta the active/selected tab instead of DatosAP
$wsCuenta.Activate()
